$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B (ID Competicao) rows 2-43 were scraped with a dropped digit (50 instead of 250).
# Restore the correct value of 250 for each data row.
$ws.Range("B2:B43").Value = 250
